$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 535.8
$ws.Range("I4").Value = 535.8
$ws.Range("K4").Value = 535.8
$ws.Range("M4").Value = -421.8
$ws.Range("H28").Value = 1701.9131
$ws.Range("I28").Value = 1270.8
$ws.Range("K28").Value = 1270.8
$ws.Range("M28").Value = -785.8
$ws.Range("H32").Value = 50001748
$ws.Range("I32").Value = 100000920
$ws.Range("K32").Value = 100000920
$ws.Range("M32").Value = -100000594
$ws.Range("H33").Value = 953.6667
$ws.Range("I33").Value = 944.4
$ws.Range("K33").Value = 944.4
$ws.Range("M33").Value = -715.4
$ws.Range("H39").Value = 38.666668
$ws.Range("I39").Value = 38.666668
$ws.Range("K39").Value = 116.000004
$ws.Range("M39").Value = 179.999996
$ws.Range("H51").Value = 5678.357
$ws.Range("I51").Value = 5408.3335
$ws.Range("J51").Value = 5880.875
$ws.Range("K51").Value = 5408.3335
$ws.Range("L51").Value = 5880.875
$ws.Range("M51").Value = -4924.3335
$ws.Range("N51").Value = -6848.875
$ws.Range("H58").Value = 2887
$ws.Range("I58").Value = 374.5
$ws.Range("J58").Value = 5399.5
$ws.Range("K58").Value = 1123.5
$ws.Range("L58").Value = 16198.5
$ws.Range("M58").Value = -973.5
$ws.Range("N58").Value = -16498.5
$ws.Range("H69").Value = 7089.8
$ws.Range("I69").Value = 5816.6665
$ws.Range("J69").Value = 8999.5
$ws.Range("K69").Value = 17449.9995
$ws.Range("L69").Value = 26998.5
$ws.Range("M69").Value = -16575.9995
$ws.Range("N69").Value = -28746.5
$ws.Range("H72").Value = 7089.8
$ws.Range("I72").Value = 5816.6665
$ws.Range("J72").Value = 8999.5
$ws.Range("K72").Value = 52349.9985
$ws.Range("L72").Value = 80995.5
$ws.Range("M72").Value = -47981.9985
$ws.Range("N72").Value = -89731.5
$ws.Range("H98").Value = 1626.0476
$ws.Range("I98").Value = 1407.35
$ws.Range("K98").Value = 1407.35
$ws.Range("M98").Value = 90.65000000000009
$ws.Range("H111").Value = 4285.933
$ws.Range("I111").Value = 1749.1428
$ws.Range("K111").Value = 5247.428400000001
$ws.Range("M111").Value = -2180.428400000001
$ws.Range("H122").Value = 1626.0476
$ws.Range("I122").Value = 1407.35
$ws.Range("K122").Value = 4222.049999999999
$ws.Range("M122").Value = -1772.049999999999
$ws.Range("H132").Value = 11859.024
$ws.Range("I132").Value = 1831.6136
$ws.Range("K132").Value = 5494.8408
$ws.Range("M132").Value = -2964.8408
$ws.Range("H138").Value = 3614.8696
$ws.Range("J138").Value = 5001.9
$ws.Range("L138").Value = 15005.7
$ws.Range("N138").Value = -25285.7
$ws.Range("H141").Value = 8213.571
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5102.5454
$ws.Range("I45").Value = 4246.2
$ws.Range("K45").Value = 4246.2
$ws.Range("M45").Value = -3869.2
$ws.Range("H61").Value = 5987.8726
$ws.Range("I61").Value = 6539.615
$ws.Range("J61").Value = 3298.125
$ws.Range("K61").Value = 6539.615
$ws.Range("L61").Value = 3298.125
$ws.Range("M61").Value = -6327.615
$ws.Range("N61").Value = -3722.125
$ws.Range("H110").Value = 1024823.9
$ws.Range("I110").Value = 1460050.9
$ws.Range("K110").Value = 1460050.9
$ws.Range("M110").Value = -1458005.9
$ws.Range("H132").Value = 11395.903
$ws.Range("I132").Value = 16912.555
$ws.Range("K132").Value = 50737.665
$ws.Range("M132").Value = -48207.665
$ws.Range("H136").Value = 5987.8726
$ws.Range("I136").Value = 6539.615
$ws.Range("J136").Value = 3298.125
$ws.Range("K136").Value = 19618.845
$ws.Range("L136").Value = 9894.375
$ws.Range("M136").Value = -17068.845
$ws.Range("N136").Value = -14994.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 574.8333
$ws.Range("I22").Value = 619.8
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 619.8
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -446.8
$ws.Range("N22").Value = -696
$ws.Range("H105").Value = 2366.8
$ws.Range("I105").Value = 2076.9285
$ws.Range("K105").Value = 2076.9285
$ws.Range("M105").Value = -329.9285
$ws.Range("H134").Value = 2269.5
$ws.Range("I134").Value = 1521.8889
$ws.Range("K134").Value = 4565.6667
$ws.Range("M134").Value = -2030.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2495.0833
$ws.Range("I16").Value = 2491.7144
$ws.Range("K16").Value = 2491.7144
$ws.Range("M16").Value = -2204.7144
$ws.Range("H113").Value = 2495.0833
$ws.Range("I113").Value = 2491.7144
$ws.Range("K113").Value = 2491.7144
$ws.Range("M113").Value = -321.7143999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 77029464
$ws.Range("J37").Value = 77029464
$ws.Range("L37").Value = 231088392
$ws.Range("N37").Value = -231088616
$ws.Range("H38").Value = 205.13333
$ws.Range("I38").Value = 104.4
$ws.Range("K38").Value = 313.2
$ws.Range("M38").Value = 33.79999999999995
$ws.Range("H86").Value = 307.5
$ws.Range("I86").Value = 315
$ws.Range("J86").Value = 300
$ws.Range("K86").Value = 945
$ws.Range("L86").Value = 900
$ws.Range("M86").Value = 241
$ws.Range("N86").Value = -3272
$ws.Range("H89").Value = 307.5
$ws.Range("I89").Value = 315
$ws.Range("J89").Value = 300
$ws.Range("K89").Value = 2835
$ws.Range("M89").Value = 3093
$ws.Range("N89").Value = -14556

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 64153.914
$ws.Range("I132").Value = 86173.52
$ws.Range("J132").Value = 9104.9
$ws.Range("K132").Value = 258520.56
$ws.Range("L132").Value = 27314.7
$ws.Range("M132").Value = -255990.56
$ws.Range("N132").Value = -32374.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 814.86664
$ws.Range("J22").Value = 806.125
$ws.Range("L22").Value = 806.125
$ws.Range("N22").Value = -1396.125
$ws.Range("H27").Value = 814.86664
$ws.Range("J27").Value = 806.125
$ws.Range("L27").Value = 806.125
$ws.Range("N27").Value = -1020.125
$ws.Range("H93").Value = 4092.2104
$ws.Range("I93").Value = 926.58826
$ws.Range("J93").Value = 31000
$ws.Range("K93").Value = 926.58826
$ws.Range("L93").Value = 31000
$ws.Range("M93").Value = 321.41174
$ws.Range("N93").Value = -33496
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("H100").Value = 5243.136
$ws.Range("I100").Value = 3594.6
$ws.Range("K100").Value = 3594.6
$ws.Range("M100").Value = -3053.6
$ws.Range("H122").Value = 35719812
$ws.Range("I122").Value = 62504524
$ws.Range("K122").Value = 187513572
$ws.Range("M122").Value = -187511122
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 100429
$ws.Range("J46").Value = 100429
$ws.Range("L46").Value = 100429
$ws.Range("N46").Value = -100891
$ws.Range("H113").Value = 835.8
$ws.Range("I113").Value = 940.7778
$ws.Range("J113").Value = 678.3333
$ws.Range("K113").Value = 2822.3334
$ws.Range("L113").Value = 2034.9999
$ws.Range("M113").Value = -652.3334
$ws.Range("N113").Value = -6374.9999
$ws.Range("H132").Value = 3061.22
$ws.Range("I132").Value = 898
$ws.Range("J132").Value = 7260.4116
$ws.Range("K132").Value = 2694
$ws.Range("L132").Value = 21781.2348
$ws.Range("M132").Value = -164
$ws.Range("N132").Value = -26841.2348
$ws.Range("H134").Value = 100429
$ws.Range("J134").Value = 100429
$ws.Range("L134").Value = 301287
$ws.Range("N134").Value = -306357
